# Apply scheduled data-runner updates to the Marilith_Profits leve-profit tables
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 200.46153
$ws.Range("I11").Value = 200.46153
$ws.Range("K11").Value = 200.46153
$ws.Range("M11").Value = -60.46153000000001

$ws.Range("H100").Value = 883
$ws.Range("I100").Value = 899.5
$ws.Range("J100").Value = 850
$ws.Range("K100").Value = 899.5
$ws.Range("L100").Value = 850
$ws.Range("M100").Value = -358.5
$ws.Range("N100").Value = -1932

$ws.Range("H111").Value = 991.5
$ws.Range("I111").Value = 991.5
$ws.Range("K111").Value = 2974.5
$ws.Range("M111").Value = 92.5

$ws.Range("H138").Value = 3246.7
$ws.Range("J138").Value = 3771.25
$ws.Range("L138").Value = 11313.75
$ws.Range("N138").Value = -21593.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 784
$ws.Range("I97").Value = 713.8333
$ws.Range("K97").Value = 713.8333
$ws.Range("M97").Value = -217.8333

$ws.Range("H121").Value = 39627.5
$ws.Range("J121").Value = 39627.5
$ws.Range("L121").Value = 39627.5
$ws.Range("N121").Value = -43121.5

$ws.Range("H122").Value = 2583.3333
$ws.Range("I122").Value = 2100
$ws.Range("K122").Value = 6300
$ws.Range("M122").Value = -3850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2877.1538
$ws.Range("I86").Value = 2940.4
$ws.Range("K86").Value = 2940.4
$ws.Range("M86").Value = -1817.4

$ws.Range("H89").Value = 2877.1538
$ws.Range("I89").Value = 2940.4
$ws.Range("K89").Value = 14702
$ws.Range("M89").Value = -9086

$ws.Range("H106").Value = 1856.6666
$ws.Range("J106").Value = 1856.6666
$ws.Range("L106").Value = 1856.6666
$ws.Range("N106").Value = -4380.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4498.9375
$ws.Range("I31").Value = 4130.25
$ws.Range("J31").Value = 5605
$ws.Range("K31").Value = 4130.25
$ws.Range("L31").Value = 5605
$ws.Range("M31").Value = -3835.25
$ws.Range("N31").Value = -6195

$ws.Range("H34").Value = 4498.9375
$ws.Range("I34").Value = 4130.25
$ws.Range("J34").Value = 5605
$ws.Range("K34").Value = 4130.25
$ws.Range("L34").Value = 5605
$ws.Range("M34").Value = -3928.25
$ws.Range("N34").Value = -6009

$ws.Range("H112").Value = 69999.5
$ws.Range("J112").Value = 69999.5
$ws.Range("L112").Value = 69999.5
$ws.Range("N112").Value = -72953.5

$ws.Range("H132").Value = 1484.1177
$ws.Range("I132").Value = 1646.4286
$ws.Range("J132").Value = 726.6667
$ws.Range("K132").Value = 4939.2858
$ws.Range("L132").Value = 2180.0001
$ws.Range("M132").Value = -2409.2858
$ws.Range("N132").Value = -7240.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 101.57895
$ws.Range("J2").Value = 144.66667
$ws.Range("L2").Value = 868.0000200000001
$ws.Range("N2").Value = -1094.00002

$ws.Range("H12").Value = 228.4
$ws.Range("J12").Value = 248.22223
$ws.Range("L12").Value = 744.66669
$ws.Range("N12").Value = -1090.66669

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H22").Value = 61230.855
$ws.Range("I22").Value = 1719.6
$ws.Range("J22").Value = 210009
$ws.Range("K22").Value = 1719.6
$ws.Range("L22").Value = 210009
$ws.Range("M22").Value = -1190.6
$ws.Range("N22").Value = -211067

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

$ws.Range("H48").Value = 29999
$ws.Range("I48").Value = 29999
$ws.Range("K48").Value = 29999
$ws.Range("M48").Value = -29514

$ws.Range("H97").Value = 679.0909
$ws.Range("I97").Value = 673.3333
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 673.3333
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -177.3333
$ws.Range("N97").Value = -1792

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 60
$ws.Range("I23").Value = 60
$ws.Range("K23").Value = 60
$ws.Range("M23").Value = 170

$ws.Range("H40").Value = 6835.727
$ws.Range("I40").Value = 7766.5
$ws.Range("J40").Value = 5718.8
$ws.Range("K40").Value = 7766.5
$ws.Range("L40").Value = 5718.8
$ws.Range("M40").Value = -7630.5
$ws.Range("N40").Value = -5990.8

$ws.Range("H68").Value = 1900
$ws.Range("I68").Value = 1900
$ws.Range("K68").Value = 1900
$ws.Range("M68").Value = -1151

$ws.Range("H71").Value = 1900
$ws.Range("I71").Value = 1900
$ws.Range("K71").Value = 9500
$ws.Range("M71").Value = -5756

$ws.Range("H100").Value = 2427
$ws.Range("I100").Value = 2598.1667
$ws.Range("K100").Value = 2598.1667
$ws.Range("M100").Value = -2057.1667

$ws.Range("H136").Value = 2176.923
$ws.Range("I136").Value = 1690.1
$ws.Range("J136").Value = 3799.6667
$ws.Range("K136").Value = 5070.299999999999
$ws.Range("L136").Value = 11399.0001
$ws.Range("M136").Value = -2520.299999999999
$ws.Range("N136").Value = -16499.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H96").Value = 1623.5
$ws.Range("I96").Value = 1249.5
$ws.Range("J96").Value = 1997.5
$ws.Range("K96").Value = 1249.5
$ws.Range("L96").Value = 1997.5
$ws.Range("M96").Value = 123.5
$ws.Range("N96").Value = -4743.5

$ws.Range("H100").Value = 998.5
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 998.5
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 1997
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -3079

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H112").Value = 45812.168
$ws.Range("J112").Value = 45812.168
$ws.Range("L112").Value = 45812.168
$ws.Range("N112").Value = -48766.168
